# Generate Report for Handback
# Update the "generate date" / handoff / handback timestamp cells that get
# refreshed each time the handback report is regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# G3: "Latest HO Xliff Generate Date" for 5fd4ca21-...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-19 04:41:09"

# --- zh-cn sheet ---
# H3: "Correspond Handoff Datetime" for 5fd4ca21-...b1b829552a103bf20367438008e70943ef54a212.zh-cn.xlf
# K3: "Correspond Handback DateTime" for the same row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-19 04:41:00"
$wsZhCn.Range("K3").Value = "2016-08-19 04:41:29"

# --- de-de sheet ---
# K3: "Correspond Handback DateTime" for 5fd4ca21-...b1b829552a103bf20367438008e70943ef54a212.de-de.xlf
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-19 04:41:36"
